# Daily update at 8 AM UTC
# Appends the new day's row (45671 -> 2025-01-14) to the "Wins Over Time"
# tracker, and shifts the "latest row" date-format highlight from the
# previous last row (83) down onto the newly added last row (84).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous last row (83) reverts to the regular "YYYY-MM-DD HH:MM:SS"
# date format used by every other data row.
$ws.Range("A83").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New row 84: next day's data.
$ws.Range("A84").Value = 45671
$ws.Range("B84").Value = 197
$ws.Range("C84").Value = 196
$ws.Range("D84").Value = 193

# The new last row gets the distinct "YYYY-MM-DD" date format.
$ws.Range("A84").NumberFormat = "YYYY-MM-DD"
